$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for columns B:F (rows 2-25)
$bfData = @(
@(1.02,1.033095748054874,1.035837824644372,1.032367888000352,1.041264067148589),
@(1.02,1.034725863497439,1.037074053980286,1.033776928272614,1.042974392707843),
@(1.02,1.035778243914784,1.037871668507611,1.034686703648658,1.044078879040008),
@(1.02,1.036220098941274,1.038206441095392,1.035068711459252,1.044542688855023),
@(1.02,1.036294255519627,1.038262619186853,1.035132825444455,1.044620534653301),
@(1.02,1.035784150209496,1.037876143883805,1.034691809857147,1.044085078500674),
@(1.02,1.03364716016147,1.036256097028968,1.032844492881393,1.041842544563812),
@(1.02,1.029862485496266,1.033383321459145,1.029573790681577,1.037873433578812),
@(1.02,1.027325815932185,1.031455489934723,1.027382312206215,1.035214815411138),
@(1.02,1.026224024956077,1.030617602143004,1.026430633345173,1.034060451929174),
@(1.02,1.025814245696935,1.030305893898078,1.026076712100225,1.033631179185868),
@(1.02,1.025902168736772,1.0303727782319,1.026152648843106,1.033723282102077),
@(1.02,1.026190163271231,1.030591846080719,1.026401386848785,1.034024978210796),
@(1.02,1.026367536192252,1.03072675718279,1.026554585833761,1.034210797649698),
@(1.02,1.027398865356491,1.031511031112012,1.02744541298648,1.035291358739342),
@(1.02,1.028044871220771,1.032002141969011,1.028003459824211,1.035968308209439),
@(1.02,1.028421348822637,1.03228829788622,1.028328694147911,1.036362857121223),
@(1.02,1.028549662988278,1.032385818951487,1.028439546053953,1.036497336958106),
@(1.02,1.027975594803011,1.031949481601221,1.027943614168881,1.035895709519933),
@(1.02,1.026105370671822,1.030527349374859,1.026328151540666,1.033936149897492),
@(1.02,1.024926441419927,1.029630419603168,1.025309979356703,1.032701249950323),
@(1.02,1.025551707736269,1.030106165831834,1.025849969400086,1.033356168810351),
@(1.02,1.028006898843596,1.03197327748399,1.027970656663211,1.035928514682806),
@(1.02,1.030843247153841,1.034128193999134,1.030421242386259,1.038901696383919)
)

# New vm_pu values for columns I:M (rows 2-25)
$imData = @(
@(1.033257053570356,1.038222354002993,1.038633284321467,1.035173321254598,1.044044074934244),
@(1.033644668577635,1.039491589352489,1.039677857957511,1.036389513817422,1.045562631427252),
@(1.033893043715278,1.040310163712611,1.040350909144308,1.037173994081587,1.046542583310784),
@(1.03399688034239,1.04065365384867,1.040633183284157,1.037503205066693,1.046953930698858),
@(1.034014281036542,1.040711290177688,1.040680538898501,1.03755844701284,1.047022961450647),
@(1.033894433460863,1.040314755940801,1.040354683553152,1.037178395299684,1.046548082192401),
@(1.033388556616968,1.038651863868987,1.038986899184738,1.035584857458017,1.04455783498056),
@(1.032478317804452,1.03570049581953,1.03655445887289,1.032757475077431,1.041029912062323),
@(1.031858628214962,1.033718117961168,1.034917404475413,1.030858985989489,1.038663221705197),
@(1.031587197210532,1.032856084945113,1.034204774506838,1.030033578471481,1.037634760636037),
@(1.031485905548131,1.032535327662818,1.033939495491647,1.02972647151333,1.037252179239909),
@(1.031507654294005,1.032604156707189,1.033996424981817,1.029792370423935,1.037334270049358),
@(1.03157883402431,1.03282958253512,1.034182858276301,1.030008203459151,1.0376031479492),
@(1.031622627786949,1.032968400300011,1.034297649304603,1.030141116903538,1.037768737186225),
@(1.031876576501218,1.033775250373485,1.034964619050841,1.03091369409555,1.038731398795248),
@(1.032035038584023,1.034280380790188,1.035381974361457,1.031397407177848,1.03933425897727),
@(1.032127167887168,1.034574663390821,1.035625047099875,1.031679226411375,1.039685544116076),
@(1.032158531016939,1.034674946888344,1.03570786724053,1.03177526511955,1.039805263815975),
@(1.032018068053971,1.034226221509051,1.035337233758533,1.03134554276884,1.03926961439939),
@(1.031557886378628,1.032763215793967,1.034127974308626,1.029944660280965,1.037523985823792),
@(1.031265829933341,1.031840121238569,1.033364324616098,1.02906089148516,1.036423162950955),
@(1.031420913993505,1.032329782540733,1.033769469559321,1.029529679953492,1.037007045311348),
@(1.03202573722976,1.0342506948428,1.035357451234523,1.031368979059058,1.03929882559223),
@(1.032715888075867,1.036466060332153,1.03718598657946,1.033490772079453,1.041944509716572)
)

$bfArr = New-Object 'object[,]' 24,5
for ($r = 0; $r -lt 24; $r++) {
    for ($c = 0; $c -lt 5; $c++) {
        $bfArr[$r, $c] = $bfData[$r][$c]
    }
}
$ws.Range("B2:F25").Value = $bfArr

$imArr = New-Object 'object[,]' 24,5
for ($r = 0; $r -lt 24; $r++) {
    for ($c = 0; $c -lt 5; $c++) {
        $imArr[$r, $c] = $imData[$r][$c]
    }
}
$ws.Range("I2:M25").Value = $imArr

Write-Output "vm_pu values updated"
